# Applies the "add 2022-Q4 data" change:
#  1. Insert a new row at the top of the "总计" (Total) sheet's data table for
#     the 2022-Q4 summary figures, shifting the existing quarters down.
#  2. Insert a brand-new worksheet named "2022-Q4" right after "总计", holding
#     the per-fund holding detail for that quarter (copied formatting from the
#     existing "2022-Q3" sheet, which sits right after it).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" summary sheet.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Make room for the new quarter as row 2 (row 1 is the header).
$total.Rows.Item(2).Insert()

# Copy the formatting of the row right below (the old row 2, now row 3) onto
# the freshly inserted row, so style (bold/border on column A, etc.) matches.
# Use a bounded A:D range (not whole-row references) to avoid ballooning the
# sheet's used range out to column XFD.
$totalSrc = $total.Range($total.Cells.Item(3, 1), $total.Cells.Item(3, 4))
$totalSrc.Copy()
$totalDst = $total.Range($total.Cells.Item(2, 1), $total.Cells.Item(2, 4))
$totalDst.PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new quarter's summary values.
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 21
$total.Cells.Item(2, 4).Value = 0.8

# Renumber the running index in column A (0-based row position) for the
# quarters that got shifted down.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(6, 1).Value = 4

# ---------------------------------------------------------------------------
# 2. Insert the new "2022-Q4" detail sheet right after "总计".
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $total)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The Q3 template only has 9 data rows (rows 2-10); Q4 needs 21 (rows 2-22).
# Extend the table by copying the format of the last template row down. Use a
# bounded A:H range (not whole-row references) to avoid ballooning the
# sheet's used range out to column XFD.
$q4Src = $q4.Range($q4.Cells.Item(10, 1), $q4.Cells.Item(10, 8))
$q4Src.Copy()
$fillRange = $q4.Range($q4.Cells.Item(11, 1), $q4.Cells.Item(22, 8))
$fillRange.PasteSpecial(-4122)  # xlPasteFormats

# Force the text columns (B-G) to keep their "text" formatting so that
# numeric-looking strings (fund codes, percentages, ...) are not coerced
# into numbers.
$textCols = $q4.Range($q4.Cells.Item(2, 2), $q4.Cells.Item(22, 7))
$textCols.NumberFormat = "@"

$q4Data = @(
    @(0,  "010434", "红土创新医疗保健股票",               "3.12", "94.35", "6.35", "0.1981", 8),
    @(1,  "014016", "中信建投品质优选一年持有期混合A",     "4.82", "76.97", "2.95", "0.1422", 2),
    @(2,  "008347", "中信建投价值甄选混合A",               "4.72", "63.36", "2.66", "0.1256", 4),
    @(3,  "003822", "中信建投行业轮换混合A",               "4.37", "56.48", "1.90", "0.0830", 9),
    @(4,  "008348", "中信建投价值甄选混合C",               "1.98", "63.36", "2.66", "0.0527", 4),
    @(5,  "003823", "中信建投行业轮换混合C",               "2.72", "56.48", "1.90", "0.0517", 9),
    @(6,  "014017", "中信建投品质优选一年持有期混合C",     "1.34", "76.97", "2.95", "0.0395", 2),
    @(7,  "002504", "鹏华金鼎灵活配置混合A",               "0.60", "88.96", "5.34", "0.0320", 3),
    @(8,  "006193", "鑫元核心资产股票A",                   "0.50", "87.44", "3.07", "0.0154", 10),
    @(9,  "007468", "中信建投策略精选混合A",               "0.58", "70.09", "2.60", "0.0151", 6),
    @(10, "007469", "中信建投策略精选混合C",               "0.41", "70.09", "2.60", "0.0107", 6),
    @(11, "002505", "鹏华金鼎灵活配置混合C",               "0.18", "88.96", "5.34", "0.0096", 3),
    @(12, "003828", "鹏华兴惠定期开放灵活配置混合",         "1.22", "20.25", "0.67", "0.0082", 4),
    @(13, "002543", "长城久益灵活配置混合A",               "0.19", "88.77", "3.21", "0.0061", 7),
    @(14, "008119", "鹏华金享混合",                       "0.76", "24.13", "0.69", "0.0052", 8),
    @(15, "002544", "长城久益灵活配置混合C",               "0.07", "88.77", "3.21", "0.0022", 7),
    @(16, "001330", "鹏华弘实灵活配置混合C",               "0.26", "25.09", "0.66", "0.0017", 8),
    @(17, "001453", "鹏华弘鑫灵活配置混合A",               "0.25", "20.45", "0.52", "0.0013", 7),
    @(18, "001454", "鹏华弘鑫灵活配置混合C",               "0.09", "20.45", "0.52", "0.0005", 7),
    @(19, "001329", "鹏华弘实灵活配置混合A",               "0.05", "25.09", "0.66", "0.0003", 8),
    @(20, "006194", "鑫元核心资产股票C",                   "0.01", "87.44", "3.07", "0.0003", 10)
)

$r = 2
foreach ($row in $q4Data) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = $row[3]
    $q4.Cells.Item($r, 5).Value = $row[4]
    $q4.Cells.Item($r, 6).Value = $row[5]
    $q4.Cells.Item($r, 7).Value = $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Restore the original active sheet/selection ("总计", cell A1) so that the
# workbook-level view state matches the source file as closely as possible.
$total.Activate()
[void]$total.Range("A1").Select()

